$p = $ppt.ActivePresentation
try {
    $p.Slides.InsertFromFile("/tmp/work/office_theme_src.pptx", 3)
    Write-Output "ok"
} catch {
    Write-Output "EXC: $_"
}
